# Fix alternate material cookbook failing with the addition of days_unit parameter
#
# - adds a new "days_unit" row to the parameter sheet
# - converts the lead-time / duration example values in "item supplier" and
#   "operation" from seconds to days (864000s -> 10d, 6048000s -> 70d, 86400s -> 1d)
# - leaves "operation" as the active sheet/tab

$wb = $excel.ActiveWorkbook

# --- parameter sheet: document the new days_unit parameter -----------------
$wsParam = $wb.Worksheets.Item("parameter")
$wsParam.Cells.Item(15, 1).Value = "days_unit"
# Copy/paste the existing "true" cell (B14) so the new cell reuses the same
# shared string as a plain text value instead of Excel auto-coercing a typed
# "true" into a boolean.
$wsParam.Range("B14").Copy()
$wsParam.Range("B15").PasteSpecial(-4163)
$wsParam.Cells.Item(15, 3).Value = "Determines whether numbers in spreadsheets are considered as days or seconds. Default is true (for days)."
$wsParam.Range("C15").Select()

# --- item supplier sheet: lead times expressed in days now -----------------
$wsItemSupplier = $wb.Worksheets.Item("item supplier")
$wsItemSupplier.Cells.Item(2, 4).Value = 10
$wsItemSupplier.Cells.Item(3, 4).Value = 10
$wsItemSupplier.Cells.Item(4, 4).Value = 10
$wsItemSupplier.Cells.Item(5, 4).Value = 70
$wsItemSupplier.Cells.Item(6, 4).Value = 10
$wsItemSupplier.Range("D5").Select()

# --- operation sheet: durations expressed in days now -----------------------
$wsOperation = $wb.Worksheets.Item("operation")
$wsOperation.Cells.Item(2, 9).Value = 1
$wsOperation.Cells.Item(3, 9).Value = 1

# Make "operation" the active sheet/tab, with I4 selected, matching the
# workbook state the fix was saved in.
$wsOperation.Activate()
$wsOperation.Range("I4").Select()
